$wb = $excel.ActiveWorkbook

# --- Update header strings on existing sheets from dot-notation to JSON Pointer notation ---

# Sheet "Missing" (sheet1)
$ws1 = $wb.Worksheets.Item("Missing")
$ws1.Range("A1").Value = "/Key"
$ws1.Range("B1").Value = "/Data/A"
$ws1.Range("C1").Value = "/AllNull"
$ws1.Range("E1").Value = "/Data/B"
$ws1.Range("E3").Select()

# Sheet "Sheet1" (sheet2)
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("A1").Value = "/a/b/c/d/e/f"
$ws2.Range("B6").Select()

# Sheet "promotion" (sheet3)
$ws3 = $wb.Worksheets.Item("promotion")
$ws3.Range("A1").Value = "/t1/A"
$ws3.Range("B1").Value = "/t1/B"
$ws3.Range("C1").Value = "/t2/A"
$ws3.Range("D1").Value = "/t2/B"
$ws3.Range("E1").Value = "/t3/A"
$ws3.Range("F1").Value = "/t3/B"
$ws3.Range("G1").Value = "/t3/C"
$ws3.Range("F5").Select()

# Sheet "mergeA" (sheet4)
$ws4 = $wb.Worksheets.Item("mergeA")
$ws4.Range("A1").Value = "/Key"
$ws4.Range("B1").Value = "/Address/State"
$ws4.Range("C1").Value = "/Address/City"
$ws4.Range("D1").Value = "/Address/TEL(Int)"
$ws4.Range("C1").Select()

# Sheet "mergeB" (sheet5)
$ws5 = $wb.Worksheets.Item("mergeB")
$ws5.Range("A1").Value = "/Key"
$ws5.Range("B1").Value = "/Name"
$ws5.Range("C1").Value = "/Property/1/A"
$ws5.Range("D1").Value = "/Property/2/A"
$ws5.PageSetup.PaperSize = 9
$ws5.PageSetup.Orientation = 1
$ws5.Range("D19").Select()

# --- Add new sheet "mergeC" (append after the last existing sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add($null, $lastSheet)
$ws6.Name = "mergeC"

$ws6.Range("A1").Value = "/Key"
$ws6.Range("B1").Value = "/Address/TEL"

$ws6.Range("A2").Value = 3
$ws6.Range("B2").Value = "010-9999-9999"

$ws6.Range("A3").Value = 1
$ws6.Range("B3").Value = "011-9999-9999"

$ws6.Range("A4").Value = 5
$ws6.Range("B4").Value = "012-9999-9999"

$ws6.Range("A5").Value = 6
$ws6.Range("B5").Value = "013-9999-9999"

$ws6.PageSetup.PaperSize = 9
$ws6.PageSetup.Orientation = 1

$ws6.Range("B5").Select()
$ws6.Activate()
